$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns for data rows remain text, matching the source formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '76.326.40'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '2.961.52'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '199.66'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").Value = '629.10'
$ws.Range("E6").Value = '  +5.43%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").Value = '0.200'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").Value = '2.960.48'
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").Value = '0.430'
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("D12").Value = '0.160'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '4.97'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("D14").Value = '3.512.17'
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").Value = '28.94'
$ws.Range("E15").Value = '  +6.00%  '
$ws.Range("D16").Value = '76.263.06'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '0.0000187'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = '2.964.34'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +4.78%  '
$ws.Range("D20").Value = '8.74'
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("D21").Value = '371.10'
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("D22").Value = '2.25'
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("D23").Value = '4.28'
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("D24").Value = '72.61'
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").Value = '3.115.10'
$ws.Range("E25").Value = '  +1.81%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '4.30'
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  -2.92%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").Value = '  +6.27%  '
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '504.06'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  +6.90%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '164.46'
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = '20.25'
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").Value = '19.97'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = '0.379'
$ws.Range("E39").Value = '  +10.55%  '
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +14.56%  '
$ws.Range("D41").Value = '183.16'
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("D42").Value = '0.110'
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '42.73'
$ws.Range("E44").Value = '  +6.52%  '
$ws.Range("D45").Value = '4.91'
$ws.Range("E45").Value = '  -1.69%  '
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("D47").Value = '1.63'
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").Value = '0.706'
$ws.Range("E48").Value = '  +7.05%  '
$ws.Range("D49").Value = '0.582'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").Value = '2.31'
$ws.Range("E50").Value = '  -0.12%  '
$ws.Range("D51").Value = '3.81'
$ws.Range("E51").Value = '  +2.61%  '
